# "break out stock.yaml completed"
# ATGL.NS 1wk stock history sheet:
#   - a handful of previously-miscounted "detect_structure"/"isPivot" flags
#     are corrected back to their steady-state values
#   - the two trailing weeks that used to be placeholders get their
#     "backup" flag finalized
#   - 26 new weekly OHLCV rows (2024-07-01 .. 2024-12-23) are appended,
#     continuing the existing table layout/columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- corrections to existing rows ------------------------------------
$ws.Range("Q55").Value = 0
$ws.Range("Q57").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("Q75").Value = 0

$ws.Range("O293").Value = 3
$ws.Range("R295").Value = 0
$ws.Range("R296").Value = 0

# ---- new weekly rows (297-322) ----------------------------------------
# columns: Datetime, Open, High, Low, Close, Adj Close, Volume, Year,
#          Month, Day, Hour, Minute, Second, Week, isPivot,
#          two_line_structure, detect_structure, backup
$newRows = @(
    @(297, 45474, 896.5499877929688, 938, 883, 891.0499877929688, $null, 7609672, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0, $null),
    @(298, 45481, 894.0499877929688, 917.4000244140625, 880.75, 893.2000122070312, $null, 6851733, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0, $null),
    @(299, 45488, 893.5, 915.7999877929688, 883.9500122070312, 888.5, $null, 11768315, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0, $null),
    @(300, 45495, 881.5499877929688, 916.8499755859375, 857, 888.4500122070312, $null, 14637362, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0, $null),
    @(301, 45502, 900, 941.9000244140625, 883, 908.9000244140625, $null, 13896579, 2024, 7, 29, 0, 0, 0, 31, 1, 0, 0, $null),
    @(302, 45509, 890, 896.4500122070312, 851, 869.8499755859375, $null, 14967918, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0, $null),
    @(303, 45516, 839.9500122070312, 867.75, 752.3499755859375, 848.1500244140625, $null, 20521989, 2024, 8, 12, 0, 0, 0, 33, 2, 0, 0, $null),
    @(304, 45523, 846.5, 868.8499755859375, 846.0499877929688, 858.2999877929688, $null, 23568301, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0, $null),
    @(305, 45530, 859.0499877929688, 865, 826, 830.3499755859375, $null, 14450961, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0, $null),
    @(306, 45537, 832, 858.9500122070312, 823.5999755859375, 826.75, $null, 12382453, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0, $null),
    @(307, 45544, 826.0499877929688, 830.4000244140625, 798, 800.75, $null, 2625880, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 2, $null),
    @(308, 45551, 802, 825, 773.0499877929688, 788.7000122070312, $null, 3511351, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0, $null),
    @(309, 45558, 845, 855, 780, 783.75, $null, 9706957, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 1, $null),
    @(310, 45565, 784.5499877929688, 794.3499755859375, 758.1500244140625, 763.4000244140625, $null, 2034822, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0, $null),
    @(311, 45572, 758.5499877929688, 773, 728.0499877929688, 757.3499755859375, $null, 2904296, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0, $null),
    @(312, 45579, 755, 763.0499877929688, 716.9500122070312, 732.5999755859375, $null, 1768320, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0, $null),
    @(313, 45586, 730, 773.7999877929688, 682.25, 724.1500244140625, $null, 6266249, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0, $null),
    @(314, 45593, 719.0999755859375, 727.7000122070312, 700.0499877929688, 722.9500122070312, $null, 2057911, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0, $null),
    @(315, 45600, 718, 752, 706.3499755859375, 712.5999755859375, $null, 1941766, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 0, $null),
    @(316, 45607, 708.75, 713.25, 673.5499877929688, 683.9000244140625, $null, 1238816, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0, $null),
    @(317, 45614, 665, 683.7000122070312, 545.75, 609.0999755859375, $null, 14114113, 2024, 11, 18, 0, 0, 0, 47, 2, 0, 0, $null),
    @(318, 45621, 636.3499755859375, 862, 572.5499877929688, 812.3499755859375, $null, 55855126, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0, $null),
    @(319, 45628, 812, 815.2999877929688, 728, 731.0499877929688, $null, 11038067, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0, $null),
    @(320, 45635, 731.0499877929688, 755.6500244140625, 708.0499877929688, 717.8499755859375, $null, 8615638, 2024, 12, 9, 0, 0, 0, 50, 0, 0, 0, $null),
    @(321, 45642, 717.75, 723.6500244140625, 666.1500244140625, 674.3499755859375, $null, 3948961, 2024, 12, 16, 0, 0, 0, 51, 0, 0, 0, $null),
    @(322, 45649, 680, 686.3499755859375, 665.25, 679.2999877929688, $null, 2797490, 2024, 12, 23, 0, 0, 0, 52, 0, 0, 0, $null)
)

$firstNewRow = 297
$lastNewRow = 322

# Match the Datetime column's existing date/time number format so the new
# rows render the same as the rest of column A.
$ws.Range("A$firstNewRow`:A$lastNewRow").NumberFormat = $ws.Range("A296").NumberFormat()

foreach ($row in $newRows) {
    $r = $row[0]
    for ($col = 1; $col -le 18; $col++) {
        $val = $row[$col]
        if ($null -eq $val) {
            $ws.Cells.Item($r, $col).Value = ""
        } else {
            $ws.Cells.Item($r, $col).Value = $val
        }
    }
}
